{"js": "// Append 27 blank paragraphs followed by a paragraph containing\n// \"Trabajo flores\" to the very end of the document body (after the\n// existing last, empty paragraph and before the section break).\nconst body = context.document.body;\n\nconst BLANK_COUNT = 27;\nfor (let i = 0; i < BLANK_COUNT; i++) {\n    body.insertParagraph(\"\", Word.InsertLocation.end);\n}\nbody.insertParagraph(\"Trabajo flores\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Append 27 blank paragraphs followed by a paragraph containing\n# \"Trabajo flores\" to the very end of the document (after the existing\n# last, empty paragraph and before the section break).\n$d = $word.ActiveDocument\n\n$blankCount = 27\nfor ($i = 0; $i -lt $blankCount; $i++) {\n    $d.Content.InsertParagraphAfter()\n}\n\n$d.Content.InsertParagraphAfter()\n$d.Content.InsertAfter(\"Trabajo flores\")\n"}
